$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 351, shifting existing rows 351-461 down to 352-462
$ws.Rows(351).Insert()

# Populate the newly inserted row 351 with the latest weekly price record
$ws.Range("A351").Value = 6
$ws.Range("B351").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C351").Value = "Metropolitana"
$ws.Range("D351").Value = 44627
$ws.Range("E351").Value = 13
$ws.Range("F351").Value = 100112044
$ws.Range("G351").Value = "Perejil"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 160
$ws.Range("K351").Value = 14000
$ws.Range("L351").Value = 15000
$ws.Range("M351").Value = 14562
$ws.Range("N351").Value = "$/docena de atados"
$ws.Range("O351").Value = "Región Metropolitana"
$ws.Range("P351").Value = 4854
$ws.Range("Q351").Value = 3
$ws.Range("R351").Value = "Hortaliza"
